$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-check pass was run, so a new timestamped column "I" is
# appended after the last existing snapshot column "H".

# Header cell I1: same look (bold/centered/bordered) as the other
# timestamp headers, carrying the new run's timestamp.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value2 = "2026-01-27 22:12:46"

# Data rows (2-100): column I records the same price last seen in
# column H (i.e. the price hasn't changed since the previous check).
$ws.Range("I2:I100").Value2 = $ws.Range("H2:H100").Value2

# Trailing rows (101-204) have no price history yet - column H is
# blank there too. Copy H's (blank) formatting into I so the new
# column has a matching blank cell for every one of those rows.
$ws.Range("H101:H204").Copy() | Out-Null
$ws.Range("I101:I204").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
